$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns:  A=SendingCluster  B=Ligand  C=Receptor  D=TargetCluster
#           E..T = Ligand-expressing cells, detection rate, avg/total expr value,
#                  avg/total derived specificity, Receptor-expressing cells, ... ,
#                  Edge avg/total expr weight, Edge avg/total derived specificity
# Each data row: row#, SendingCluster, TargetCluster, E, F, G, H, I, J, K, L, M, N, O, P, Q, R, S, T
$rows = @(
  ,@(2, "ECs", "M1", 3, 1, 2.575111333333334, 7.725334, 0.002862529524196255, 0.002862529524196256, 2, 0.6666666666666666, 0.8058999999999999, 2.4177, 0.09451456273674329, 0.0945145627367433, 2.075282223533333, 18.6775400118, 0.0002705507263004269, 0.000270550726300427)
  ,@(3, "ECs", "M2", 3, 1, 2.575111333333334, 7.725334, 0.002862529524196255, 0.002862529524196256, 3, 1, 2.375759, 7.127277, 0.2786249200308754, 0.2786249200308755, 6.117843926168667, 55.060595335518, 0.0007975720597652014, 0.0007975720597652018)
  ,@(4, "ECs", "Neutro", 3, 1, 2.575111333333334, 7.725334, 0.002862529524196255, 0.002862529524196256, 3, 1, 5.345069333333334, 16.035208, 0.6268605172323812, 0.6268605172323812, 13.76414861771911, 123.877337559472, 0.001794406738130627, 0.001794406738130627)
  ,@(5, "FAPs", "M1", 3, 1, 378.538676, 1135.616028, 0.4207888498154877, 0.4207888498154878, 2, 0.6666666666666666, 0.8058999999999999, 2.4177, 0.09451456273674329, 0.0945145627367433, 305.0643189883999, 2745.5788708956, 0.03977067414480796, 0.03977067414480797)
  ,@(6, "FAPs", "M2", 3, 1, 378.538676, 1135.616028, 0.4207888498154877, 0.4207888498154878, 3, 1, 2.375759, 7.127277, 0.2786249200308754, 0.2786249200308755, 899.3166663550838, 8093.849997195755, 0.1172422596297243, 0.1172422596297243)
  ,@(7, "FAPs", "Neutro", 3, 1, 378.538676, 1135.616028, 0.4207888498154877, 0.4207888498154878, 3, 1, 5.345069333333334, 16.035208, 0.6268605172323812, 0.6268605172323812, 2023.315468568203, 18209.83921711382, 0.2637759160409554, 0.2637759160409555)
  ,@(8, "M1", "M1", 3, 1, 280.8555153333334, 842.566546, 0.312202891684044, 0.312202891684044, 2, 0.6666666666666666, 0.8058999999999999, 2.4177, 0.09451456273674329, 0.0945145627367433, 226.3414598071333, 2037.0731382642, 0.02950771979266425, 0.02950771979266425)
  ,@(9, "M1", "M2", 3, 1, 280.8555153333334, 842.566546, 0.312202891684044, 0.312202891684044, 3, 1, 2.375759, 7.127277, 0.2786249200308754, 0.2786249200308755, 667.2450182528047, 6005.205164275242, 0.08698750572887483, 0.08698750572887484)
  ,@(10, "M1", "Neutro", 3, 1, 280.8555153333334, 842.566546, 0.312202891684044, 0.312202891684044, 3, 1, 5.345069333333334, 16.035208, 0.6268605172323812, 0.6268605172323812, 1501.19220210573, 13510.72981895157, 0.1957076661625049, 0.1957076661625049)
  ,@(11, "M2", "M1", 3, 1, 98.90178066666668, 296.705342, 0.1099405930490186, 0.1099405930490186, 2, 0.6666666666666666, 0.8058999999999999, 2.4177, 0.09451456273674329, 0.0945145627367433, 79.70494503926668, 717.3445053534001, 0.01039098707904624, 0.01039098707904624)
  ,@(12, "M2", "M2", 3, 1, 98.90178066666668, 296.705342, 0.1099405930490186, 0.1099405930490186, 3, 1, 2.375759, 7.127277, 0.2786249200308754, 0.2786249200308755, 234.9667955348594, 2114.701159813734, 0.03063218894642983, 0.03063218894642984)
  ,@(13, "M2", "Neutro", 3, 1, 98.90178066666668, 296.705342, 0.1099405930490186, 0.1099405930490186, 3, 1, 5.345069333333334, 16.035208, 0.6268605172323812, 0.6268605172323812, 528.6368748534597, 4757.731873681137, 0.06891741702354285, 0.06891741702354257)
  ,@(14, "Neutro", "M1", 3, 1, 132.0794193333333, 396.238258, 0.1468213169995101, 0.1468213169995101, 2, 0.6666666666666666, 0.8058999999999999, 2.4177, 0.09451456273674329, 0.0945145627367433, 106.4428040407333, 957.9852363666, 0.01387675257664147, 0.01387675257664147)
  ,@(15, "Neutro", "M2", 3, 1, 132.0794193333333, 396.238258, 0.1468213169995101, 0.1468213169995101, 3, 1, 2.375759, 7.127277, 0.2786249200308754, 0.2786249200308755, 313.7888691959407, 2824.099822763466, 0.0409080777078163, 0.04090807770781631)
  ,@(16, "Neutro", "Neutro", 3, 1, 132.0794193333333, 396.238258, 0.1468213169995101, 0.1468213169995101, 3, 1, 5.345069333333334, 16.035208, 0.6268605172323812, 0.6268605172323812, 705.9736538430739, 6353.762884587665, 0.09203648671505228, 0.09203648671505231)
  ,@(17, "sCs", "M1", 3, 1, 6.642431333333334, 19.927294, 0.007383818927743304, 0.007383818927743305, 2, 0.6666666666666666, 0.8058999999999999, 2.4177, 0.09451456273674329, 0.0945145627367433, 5.353135411533334, 48.1782187038, 0.0006978784172829471, 0.0006978784172829473)
  ,@(18, "sCs", "M2", 3, 1, 6.642431333333334, 19.927294, 0.007383818927743304, 0.007383818927743305, 3, 1, 2.375759, 7.127277, 0.2786249200308754, 0.2786249200308755, 15.78081602204867, 142.027344198438, 0.002057315958264942, 0.002057315958264943)
  ,@(19, "sCs", "Neutro", 3, 1, 6.642431333333334, 19.927294, 0.007383818927743304, 0.007383818927743305, 3, 1, 5.345069333333334, 16.035208, 0.6268605172323812, 0.6268605172323812, 35.50425601857245, 319.538304167152, 0.004628624552195414, 0.004628624552195415)
)

foreach ($r in $rows) {
    $rowNum = $r[0]
    $ws.Cells.Item($rowNum, 1).Value = $r[1]   # A: Sending cluster
    $ws.Cells.Item($rowNum, 2).Value = "Ccl2"  # B: Ligand symbol
    $ws.Cells.Item($rowNum, 3).Value = "Cxcr3" # C: Receptor symbol
    $ws.Cells.Item($rowNum, 4).Value = $r[2]   # D: Target cluster

    for ($i = 0; $i -lt 16; $i++) {
        $ws.Cells.Item($rowNum, 5 + $i).Value = $r[3 + $i]
    }
}
